# Applies the typo-fix / proofErr-wrap edits described by the diff.
#
# The documented Word object model has no call that authors a <w:proofErr/>
# marker directly (those are Word's own "this span has already been
# spell-checked" bookkeeping, produced internally by the spell checker -
# not something a script inserts through a property or method). What *is*
# available is Range.InsertXML, which accepts raw WordprocessingML. When a
# Range currently spans some text, clearing that text and then calling
# InsertXML with a full "<w:p>...</w:p>" replacement swaps that paragraph's
# content in place (same slot among its siblings) without disturbing
# neighboring paragraphs. We use that to rebuild each affected paragraph
# with the exact run / proofErr structure the diff calls for, re-using the
# paragraph's original attributes so nothing else about it changes.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($index, $innerXml) {
    $p = $d.Paragraphs($index)
    $r = $p.Range
    $r.Text = ""
    $r.InsertXML($innerXml)
}

# Paragraph 10: "Инкриментная модель" -> "Инкрементная" + " модель"
# (typo fix, split across two runs; pPr/rPr unchanged)
$xml10 = '<w:p ' + $wNs + ' w:rsidR="0082367F" w:rsidRDefault="0082367F" w:rsidP="0082367F">' +
            '<w:pPr>' +
                '<w:pStyle w:val="a3"/>' +
                '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
                '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
            '</w:pPr>' +
            '<w:r><w:t>Инкрементная</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> модель</w:t></w:r>' +
         '</w:p>'
Set-ParagraphXml 10 $xml10

# Paragraph 13: "Хорошо подходит под любой обьём работ" -> split into 3 runs,
# with "обьём" wrapped between spellStart/spellEnd proofErr markers.
$xml13 = '<w:p ' + $wNs + ' w:rsidR="0082367F" w:rsidRDefault="0082367F" w:rsidP="0082367F">' +
            '<w:pPr>' +
                '<w:pStyle w:val="a5"/>' +
                '<w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr>' +
            '</w:pPr>' +
            '<w:r><w:t xml:space="preserve">Хорошо подходит под любой </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>обьём</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> работ</w:t></w:r>' +
         '</w:p>'
Set-ParagraphXml 13 $xml13

# Paragraph 19: the page-break paragraph loses its
# <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>.
$xml19 = '<w:p ' + $wNs + ' w:rsidR="0082367F" w:rsidRPr="004C1F96" w:rsidRDefault="004C1F96" w:rsidP="004C1F96">' +
            '<w:r><w:br w:type="page"/></w:r>' +
         '</w:p>'
Set-ParagraphXml 19 $xml19

# Paragraph 24: "Оринетирована на тестирование и бесперебойное функционированние"
# -> split into 3 runs, with "Оринетирована" and "функционированние" each
# wrapped between spellStart/spellEnd proofErr markers.
$xml24 = '<w:p ' + $wNs + ' w:rsidR="0082367F" w:rsidRDefault="0082367F" w:rsidP="0082367F">' +
            '<w:pPr>' +
                '<w:pStyle w:val="a5"/>' +
                '<w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr>' +
            '</w:pPr>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>Оринетирована</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> на тестирование и бесперебойное </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>функционированние</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
         '</w:p>'
Set-ParagraphXml 24 $xml24

# Paragraph 29: "Выбрана инриментная модель" -> split into 3 runs, with
# "инриментная" wrapped between spellStart/spellEnd proofErr markers.
$xml29 = '<w:p ' + $wNs + ' w:rsidR="0082367F" w:rsidRPr="004C1F96" w:rsidRDefault="0082367F" w:rsidP="0082367F">' +
            '<w:pPr>' +
                '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
            '</w:pPr>' +
            '<w:r><w:t xml:space="preserve">Выбрана </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>инриментная</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> модель</w:t></w:r>' +
         '</w:p>'
Set-ParagraphXml 29 $xml29
